# StagingTemplates/Staging.OutputOutputLink.xlsx
#
# The commit relabels the two "key" rows on Sheet1:
#   A2 used to read "OutputOutputLink_ID"  -> now reads "OutputBusinessKey"
#   B2 used to read "OutputSourceKey"      -> now reads "OutputOutputLink_ID"
#
# (Everything else in the upstream diff - the bookViews window size, the
# sheet's internal VBA codeName, and dropping column B's explicit width -
# are cosmetic artifacts of Excel re-saving the file and carry no
# spreadsheet-visible effect, so the cell content below is the
# substantive change.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "OutputBusinessKey"
$ws.Range("B2").Value = "OutputOutputLink_ID"
